$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that are being "commented out" (removed):
#  - row 7: YahooFinanceUrl
#  - row 11: QuoteUrlTemplate
#  - row 12: QuoteSummaryUrlTemplate
# Deleting row 12 first, then 11, then 7 keeps the remaining row numbers stable
# while each delete is performed.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(7).Delete()

# Update the selection to match the post-edit state (A7:XFD7, i.e. whole row 7 selected)
$ws.Range("A7:XFD7").Select()
